{"js": "// Apply the three writeup edits:\n//  1) Objective paragraph: merge the \"support or disprove...\" runs into one\n//     (no visible text change, just tidies up the grammar-check split runs).\n//  2) Race paragraph: \"were White.\" -> \"were attributed to White race.\"\n//     and \"White population is 80% for Connecticut\" -> \"the White\n//     population of Connecticut is 80%\".\n//  3) Age Groups paragraph: \"This group can\" -> \"These groups can\",\n//     \"three drugs\" -> \"top three drugs\", and \"the most deaths are\n//     illegal, it is\" -> \"the accidental drug-related deaths are\n//     illegal,  it is\" (note the extra space before \"it\").\n\nasync function replaceOnce(body, searchText, replacementText, searchOptions) {\n  const results = body.search(searchText, searchOptions || { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + searchText);\n  }\n  results.items[0].insertText(replacementText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nconst body = context.document.body;\n\n// 1) Objective paragraph - collapse the split runs back into one run with\n// the identical combined text (removes the old grammar-check run breaks).\nawait replaceOnce(\n  body,\n  \" support or disprove initial hypotheses, and answer questions about the reported deaths. \",\n  \" support or disprove initial hypotheses, and answer questions about the reported deaths. \"\n);\n\n// 2a) Race paragraph - \"were White.\" -> \"were attributed to White race.\"\nawait replaceOnce(\n  body,\n  \"78% of deaths were White.\",\n  \"78% of deaths were attributed to White race.\"\n);\n\n// 2b) Race paragraph - reword the census sentence.\nawait replaceOnce(\n  body,\n  \"White population is 80% for Connecticut\",\n  \"the White population of Connecticut is 80%\"\n);\n\n// 3a) Age Groups paragraph - \"This group can\" -> \"These groups can\"\nawait replaceOnce(\n  body,\n  \"age group.  This group can be expected\",\n  \"age group.  These groups can be expected\"\n);\n\n// 3b) Age Groups paragraph - \"three drugs\" -> \"top three drugs\"\nawait replaceOnce(\n  body,\n  \"Since two of the three drugs involved in the\",\n  \"Since two of the top three drugs involved in the\"\n);\n\n// 3c) Age Groups paragraph - \"the most deaths are illegal, it is\" ->\n// \"the accidental drug-related deaths are illegal,  it is\"\nawait replaceOnce(\n  body,\n  \"involved in the most deaths are illegal, it is reasonable\",\n  \"involved in the accidental drug-related deaths are illegal,  it is reasonable\"\n);\n", "ps1": "# Apply the three writeup edits:\n#  1) Objective paragraph: merge the \"support or disprove...\" runs into one\n#     (no visible text change, just tidies up the grammar-check split runs).\n#  2) Race paragraph: \"were White.\" -> \"were attributed to White race.\"\n#     and \"White population is 80% for Connecticut\" -> \"the White\n#     population of Connecticut is 80%\".\n#  3) Age Groups paragraph: \"This group can\" -> \"These groups can\",\n#     \"three drugs\" -> \"top three drugs\", and \"the most deaths are\n#     illegal, it is\" -> \"the accidental drug-related deaths are\n#     illegal,  it is\" (note the extra space before \"it\").\n\n$wdFindContinue = 1\n$wdReplaceOne = 1\n$wdReplaceAll = 2\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $range = $d.Content\n    $ok = $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $replaceText, $wdReplaceAll)\n    if (-not $ok) {\n        throw \"Text not found: $findText\"\n    }\n}\n\n# 1) Objective paragraph - collapse the split runs back into one run with\n# the identical combined text (removes the old grammar-check run breaks).\nReplace-Text \" support or disprove initial hypotheses, and answer questions about the reported deaths. \" \" support or disprove initial hypotheses, and answer questions about the reported deaths. \"\n\n# 2a) Race paragraph - \"were White.\" -> \"were attributed to White race.\"\nReplace-Text \"78% of deaths were White.\" \"78% of deaths were attributed to White race.\"\n\n# 2b) Race paragraph - reword the census sentence.\nReplace-Text \"White population is 80% for Connecticut\" \"the White population of Connecticut is 80%\"\n\n# 3a) Age Groups paragraph - \"This group can\" -> \"These groups can\"\nReplace-Text \"age group.  This group can be expected\" \"age group.  These groups can be expected\"\n\n# 3b) Age Groups paragraph - \"three drugs\" -> \"top three drugs\"\nReplace-Text \"Since two of the three drugs involved in the\" \"Since two of the top three drugs involved in the\"\n\n# 3c) Age Groups paragraph - \"the most deaths are illegal, it is\" ->\n# \"the accidental drug-related deaths are illegal,  it is\"\nReplace-Text \"involved in the most deaths are illegal, it is reasonable\" \"involved in the accidental drug-related deaths are illegal,  it is reasonable\"\n"}
